$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174
$ws.Range("A174").Value = 173.0
$ws.Range("B174").Value = "Sunday, Jan 15"
$ws.Range("C174").Value = "9:15 AM"
$ws.Range("D174").Value = "FR1951"
$ws.Range("E174").Value = "Manchester"
$ws.Range("F174").Value = "(MAN)"
$ws.Range("G174").Value = "Ryanair "
$ws.Range("H174").Value = "B738"
$ws.Range("I174").Value = "(EI-EBW)"
$ws.Range("J174").Value = "9:01 AM"
$ws.Range("L174").Value = "0 hours, -14 minutes"

# Row 175
$ws.Range("A175").Value = 174.0
$ws.Range("B175").Value = "Sunday, Jan 15"
$ws.Range("C175").Value = "9:35 AM"
$ws.Range("D175").Value = "FR2019"
$ws.Range("E175").Value = "London"
$ws.Range("F175").Value = "(STN)"
$ws.Range("G175").Value = "Ryanair "
$ws.Range("H175").Value = "B38M"
$ws.Range("I175").Value = "(EI-HAX)"
$ws.Range("J175").Value = "9:17 AM"
$ws.Range("L175").Value = "0 hours, -18 minutes"

# Row 176
$ws.Range("A176").Value = 175.0
$ws.Range("B176").Value = "Sunday, Jan 15"
$ws.Range("C176").Value = "11:05 AM"
$ws.Range("D176").Value = "FR8083"
$ws.Range("E176").Value = "Birmingham"
$ws.Range("F176").Value = "(BHX)"
$ws.Range("G176").Value = "Ryanair "
$ws.Range("H176").Value = "B738"
$ws.Range("I176").Value = "(EI-EXD)"
$ws.Range("J176").Value = "10:51 AM"
$ws.Range("L176").Value = "0 hours, -14 minutes"

# Row 177
$ws.Range("A177").Value = 176.0
$ws.Range("B177").Value = "Sunday, Jan 15"
$ws.Range("C177").Value = "11:50 AM"
$ws.Range("D177").Value = "FR4529"
$ws.Range("E177").Value = "Oslo"
$ws.Range("F177").Value = "(TRF)"
$ws.Range("G177").Value = "Ryanair "
$ws.Range("H177").Value = "B738"
$ws.Range("I177").Value = "(SP-RKP)"
$ws.Range("J177").Value = "11:57 AM"
$ws.Range("L177").Value = "0 hours, 7 minutes"

# Row 178
$ws.Range("A178").Value = 177.0
$ws.Range("B178").Value = "Sunday, Jan 15"
$ws.Range("C178").Value = "12:30 PM"
$ws.Range("D178").Value = "FR2027"
$ws.Range("E178").Value = "Dublin"
$ws.Range("F178").Value = "(DUB)"
$ws.Range("G178").Value = "Buzz "
$ws.Range("H178").Value = "B38M"
$ws.Range("I178").Value = "(SP-RZG)"
$ws.Range("J178").Value = "12:35 PM"
$ws.Range("L178").Value = "0 hours, 5 minutes"

# Row 179
$ws.Range("A179").Value = 178.0
$ws.Range("B179").Value = "Sunday, Jan 15"
$ws.Range("C179").Value = "1:55 PM"
$ws.Range("D179").Value = "FR2351"
$ws.Range("E179").Value = "Shannon"
$ws.Range("F179").Value = "(SNN)"
$ws.Range("G179").Value = "Ryanair "
$ws.Range("H179").Value = "B38M"
$ws.Range("I179").Value = "(SP-RZO)"
$ws.Range("J179").Value = "1:29 PM"
$ws.Range("L179").Value = "0 hours, -26 minutes"

# Row 180
$ws.Range("A180").Value = 179.0
$ws.Range("B180").Value = "Sunday, Jan 15"
$ws.Range("C180").Value = "2:05 PM"
$ws.Range("D180").Value = "FR9258"
$ws.Range("E180").Value = "Malta"
$ws.Range("F180").Value = "(MLA)"
$ws.Range("G180").Value = "Ryanair "
$ws.Range("H180").Value = "B738"
$ws.Range("I180").Value = "(SP-RSP)"
$ws.Range("J180").Value = "1:59 PM"
$ws.Range("L180").Value = "0 hours, -6 minutes"

# Row 181
$ws.Range("A181").Value = 180.0
$ws.Range("B181").Value = "Sunday, Jan 15"
$ws.Range("C181").Value = "2:15 PM"
$ws.Range("D181").Value = "FR1943"
$ws.Range("E181").Value = "Bologna"
$ws.Range("F181").Value = "(BLQ)"
$ws.Range("G181").Value = "Ryanair "
$ws.Range("H181").Value = "B738"
$ws.Range("I181").Value = "(9H-QDA)"
$ws.Range("J181").Value = "1:49 PM"
$ws.Range("L181").Value = "0 hours, -26 minutes"

# Row 182
$ws.Range("A182").Value = 181.0
$ws.Range("B182").Value = "Sunday, Jan 15"
$ws.Range("C182").Value = "2:45 PM"
$ws.Range("D182").Value = "FR1895"
$ws.Range("E182").Value = "Amman"
$ws.Range("F182").Value = "(AMM)"
$ws.Range("G182").Value = "Ryanair "
$ws.Range("H182").Value = "B738"
$ws.Range("I182").Value = "(SP-RSV)"
$ws.Range("J182").Value = "2:37 PM"
$ws.Range("L182").Value = "0 hours, -8 minutes"
